# Update "想去人数" (F column) figures for several events that appear on
# both the "展览" sheet and the combined "全部类型" sheet.
#
# Sheet "展览" (sheet1):
#   F2: 81    -> 82
#   F3: 11856 -> 11862
#   F4: 4     -> 6
#   F8: 11772 -> 11784
#   F10: 1173 -> 1172
#   F12: 51   -> 53
#   F14: 5834 -> 5837
#
# Sheet "全部类型" (sheet4):
#   F3: 81    -> 82
#   F5: 11856 -> 11862
#   F6: 4     -> 6
#   F11: 11772 -> 11784
#   F13: 1173 -> 1172
#   F15: 51   -> 53
#   F18: 5834 -> 5837

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 82
$wsExhibition.Range("F3").Value = 11862
$wsExhibition.Range("F4").Value = 6
$wsExhibition.Range("F8").Value = 11784
$wsExhibition.Range("F10").Value = 1172
$wsExhibition.Range("F12").Value = 53
$wsExhibition.Range("F14").Value = 5837

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F3").Value = 82
$wsAllTypes.Range("F5").Value = 11862
$wsAllTypes.Range("F6").Value = 6
$wsAllTypes.Range("F11").Value = 11784
$wsAllTypes.Range("F13").Value = 1172
$wsAllTypes.Range("F15").Value = 53
$wsAllTypes.Range("F18").Value = 5837
